$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Nitrate
$ws.Range("A2").Value = "Nitrate_c__Day_sp_exchange"
$ws.Range("B2").Value = -0.1375282857142793
$ws.Range("C2").Value = -0.1401839999999979
$ws.Range("D2").Value = -0.04086264285660046

# Row 3 - Glycine
$ws.Range("A3").Value = "Glycine_c__Day_sp_exchange"
$ws.Range("B3").Value = -0.01872199999999186
$ws.Range("C3").Value = -0.01872200000000252
$ws.Range("D3").Value = 0.01872199999999683

# Row 4 - L-Alanine
$ws.Range("A4").Value = "L-Alanine_c__Day_sp_exchange"
$ws.Range("B4").Value = -1.633878365079167
$ws.Range("C4").Value = -3.725505666666646
$ws.Range("D4").Value = -1.438243583333793

# Row 5 - D-Fructose
$ws.Range("A5").Value = "D-Fructose__Day_sp_exchange"
$ws.Range("B5").Value = -0.033165
$ws.Range("C5").Value = -0.5590200000000399
$ws.Range("D5").Value = 0.5590199999981873

# Row 6 - Citrate
$ws.Range("A6").Value = "Citrate_c__Day_sp_exchange"
$ws.Range("B6").Value = -1.38879330158621
$ws.Range("C6").Value = -7.734369500004439
$ws.Range("D6").Value = -1.345983222221735
